$d = $word.ActiveDocument

# 1) Remove the old "_GoBack" bookmark; we'll re-add it later at its new location.
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}

# 2) "O QUE FOI EXECUTADO FORA DO PLANEJAMENTO:" block: replace the placeholder
#    "---/---/---" lines with the real status lines. Scope the Find to the
#    text right after the (bold) header so the new text does not inherit the
#    header's bold formatting.
$headerRange1 = $d.Content
$headerRange1.Find.Execute("O QUE FOI EXECUTADO FORA DO PLANEJAMENTO:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$body1 = $d.Range($headerRange1.End, $d.Content.End)
$body1.Find.Execute(
    "---^l---^l---",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Validação Modelagem lógica^l- Inicio do Site Estático Dashboard (Gráfico com ChartJS) - Local",
    1) | Out-Null

# 3) "PLANEJAMENTO DA SEMANA 10/10 – 17/10:" block: replace the placeholder
#    "---/---/---" lines with the real plan lines for next week. Same
#    header-skipping trick as above.
$headerRange2 = $d.Content
$headerRange2.Find.Execute("PLANEJAMENTO DA SEMANA 10/10 – 17/10:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$body2 = $d.Range($headerRange2.End, $d.Content.End)
$body2.Find.Execute(
    "---^l---^l---",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Especificação da Dashboards - (Responsável: Daniel Sena)^l- Finalizar o Site Estático Dashboard (Gráfico com ChartJS) – Local - (Responsável: Gabriel Lima)^l- Validar a solução técnica - (Responsável: Erick Lee)^l- Métricas KPI´S - (Responsável: Erick Lee)^l- Iniciar criação da apresentação PPTX - (Responsável: Gustavo Castro)^l- Validar documentação - (Responsável: Daniel Sena)^l- Finalização do site estático (Simulador Financeiro) – (Responsável: Leandro Boneto)^l- Organização das ferramentas de gestão - (Responsável: Everton Barbosa)",
    1) | Out-Null

# 4) Merge away the trailing empty paragraph (it disappears in the target doc).
$lastIdx = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($lastIdx - 1)
$mark = $d.Range($secondLast.Range.End - 1, $secondLast.Range.End)
$mark.Delete() | Out-Null

# 5) Re-add the "_GoBack" bookmark at its new spot: right before "(Responsável:" in
#    the "Iniciar criação da apresentação PPTX" line.
$r2 = $d.Content
$r2.Find.Execute("- Iniciar criação da apresentação PPTX - ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($r2.End, $r2.End)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
